# Update country data rows based on the new COVID-19 statistics snapshot.
# Also renames a few country labels where the source re-sorted rows by total cases,
# causing three pairs of adjacent countries to swap positions (Indonesia/Filipinas/
# Emiratos Arabes Unidos, Finlandia/Luxemburgo, Kenia/Sri Lanka, Liberia/Islas Caimanes).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 13
$ws.Range("B13").Value = 36138
$ws.Range("C13").Value = 1329
$ws.Range("D13").Value = 7961
$ws.Range("E13").Value = 23014
$ws.Range("F13").Value = 1140
$ws.Range("G13").Value = 306
$ws.Range("H13").Value = 5163

# Row 20
$ws.Range("B20").Value = 14510
$ws.Range("C20").Value = 34
$ws.Range("E20").Value = 5114

# Row 39
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 5923
$ws.Range("C39").Value = 407
$ws.Range("D39").Value = 607
$ws.Range("E39").Value = 4796
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 24
$ws.Range("H39").Value = 520

# Row 40
$ws.Range("A40").Value = "Filipinas"
$ws.Range("B40").Value = 5878
$ws.Range("C40").Value = 218
$ws.Range("D40").Value = 487
$ws.Range("E40").Value = 5004
$ws.Range("G40").Value = 25
$ws.Range("H40").Value = 387

# Row 41
$ws.Range("A41").Value = "Emiratos Arabes Unidos"
$ws.Range("B41").Value = 5825
$ws.Range("D41").Value = 1095
$ws.Range("E41").Value = 4695
$ws.Range("F41").Value = 1
$ws.Range("H41").Value = 35

# Row 50
$ws.Range("A50").Value = "Finlandia"
$ws.Range("B50").Value = 3489
$ws.Range("C50").Value = 120
$ws.Range("D50").Value = 1700
$ws.Range("E50").Value = 1714
$ws.Range("F50").Value = 76
$ws.Range("H50").Value = 75

# Row 51
$ws.Range("A51").Value = "Luxemburgo"
$ws.Range("B51").Value = 3444
$ws.Range("D51").Value = 552
$ws.Range("E51").Value = 2823
$ws.Range("F51").Value = 35
$ws.Range("H51").Value = 69

# Row 57
$ws.Range("B57").Value = 2528
$ws.Range("C57").Value = 245
$ws.Range("D57").Value = 273
$ws.Range("E57").Value = 2122
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 133

# Row 79
$ws.Range("B79").Value = 1022
$ws.Range("C79").Value = 4
$ws.Range("D79").Value = 533
$ws.Range("E79").Value = 485

# Row 97
$ws.Range("B97").Value = 539
$ws.Range("C97").Value = 21
$ws.Range("D97").Value = 283
$ws.Range("E97").Value = 230

# Row 116
$ws.Range("D116").Value = 198
$ws.Range("E116").Value = 70

# Row 117
$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 246
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 53
$ws.Range("E117").Value = 182
$ws.Range("F117").Value = 2
$ws.Range("H117").Value = 11

# Row 118
$ws.Range("A118").Value = "Sri Lanka"
$ws.Range("B118").Value = 238
$ws.Range("D118").Value = 68
$ws.Range("E118").Value = 163
$ws.Range("F118").Value = 1
$ws.Range("H118").Value = 7

# Row 143
$ws.Range("B143").Value = 83
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 35
$ws.Range("F143").Value = 9

# Row 148
$ws.Range("A148").Value = "Liberia"
$ws.Range("B148").Value = 73
$ws.Range("C148").Value = 14
$ws.Range("D148").Value = 4
$ws.Range("E148").Value = 63
$ws.Range("F148").Value = 0
$ws.Range("H148").Value = 6

# Row 149
$ws.Range("A149").Value = "Islas Caimanes"
$ws.Range("B149").Value = 61
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 7
$ws.Range("E149").Value = 53
$ws.Range("F149").Value = 3
$ws.Range("H149").Value = 1

